$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks numeric need NumberFormat forced to Text
# so Excel stores them as strings (preserving exact formatting, e.g.
# trailing zeros) instead of auto-converting to a number; the style is
# then reset back to Normal so no stray style index is introduced.

$ws.Range('D2').Value = '64.949.57'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '3.146.54'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.25'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.67'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.28%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.145.59'
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('E10').Value = '  -2.90%  '
$ws.Range('E11').Value = '  -1.00%  '
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000264'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.08'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.46%  '
$ws.Range('D15').Value = '3.661.30'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('D16').Value = '64.894.41'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').Value = '3.150.74'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '504.15'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.04'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('E22').Value = '  -3.55%  '
$ws.Range('E23').Value = '  -2.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.72'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.18'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.06'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.80%  '
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('E29').Value = '  -1.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.78'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.50'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.35'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.48'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '54.90'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.76%  '
$ws.Range('E37').Value = '  +1.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '474.30'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.05%  '
$ws.Range('E39').Value = '  -2.80%  '
$ws.Range('E40').Value = '  -3.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.73'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.44%  '
$ws.Range('D42').Value = '2.991.40'
$ws.Range('E42').Value = '  -4.29%  '
$ws.Range('E43').Value = '  -2.85%  '
$ws.Range('E44').Value = '  -4.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.42'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.22'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -5.16%  '
$ws.Range('D47').Value = '0.0₃0591'
$ws.Range('E47').Value = '  +1.47%  '
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('E49').Value = '  -1.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.25'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.11%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '118.70'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.77%  '
